# Generate Report for Handback
# ------------------------------------------------------------------
# This script updates the localization-status workbook to reflect
# that the zh-cn and de-de handback packages have been produced:
#   * Status columns flip from "Ready for handoff" to
#     "Handed back: in sync with en-US" (Overview + per-language sheets)
#   * The per-language sheets grow two new populated columns:
#       F = Latest Target File   (the .md source, now hyperlinked)
#       G = Latest Handback File (the translated .xlf, now hyperlinked)
#   * Latest Handback DateTime (column H) is stamped with the real
#     handback timestamp instead of the "0001-01-01 00:00:00" placeholder.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------
# Overview sheet: flip the per-language status cells
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("H2").Value = "2016-03-14 04:30:29"
$wsZh.Range("H3").Value = "2016-03-14 04:30:29"

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/d641bdf9e317d10099a875bca9a537f3b06a36e7/e2e/8e2b4336-5735-4c5d-bffc-e21e7f60d885.md", [Type]::Missing, [Type]::Missing, "8e2b4336-5735-4c5d-bffc-e21e7f60d885.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dd5a2f21995a0c53a3e398d6721c8416c3088b66/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8e2b4336-5735-4c5d-bffc-e21e7f60d885.2a5251d41fd3c63e7892720e50d02aaef985fb9a.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "8e2b4336-5735-4c5d-bffc-e21e7f60d885.2a5251d41fd3c63e7892720e50d02aaef985fb9a.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/d641bdf9e317d10099a875bca9a537f3b06a36e7/e2e/e0390748-aa6c-4338-96c4-c0d594e4314f.md", [Type]::Missing, [Type]::Missing, "e0390748-aa6c-4338-96c4-c0d594e4314f.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dd5a2f21995a0c53a3e398d6721c8416c3088b66/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e0390748-aa6c-4338-96c4-c0d594e4314f.6f38887ef8fb072b659eafcafeb8544e61d5db31.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "e0390748-aa6c-4338-96c4-c0d594e4314f.6f38887ef8fb072b659eafcafeb8544e61d5db31.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("H2").Value = "2016-03-14 04:30:34"
$wsDe.Range("H3").Value = "2016-03-14 04:30:34"

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/d641bdf9e317d10099a875bca9a537f3b06a36e7/e2e/8e2b4336-5735-4c5d-bffc-e21e7f60d885.md", [Type]::Missing, [Type]::Missing, "8e2b4336-5735-4c5d-bffc-e21e7f60d885.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aa4e0fd75e94cdda9e4e2524ee477d04f655959d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8e2b4336-5735-4c5d-bffc-e21e7f60d885.2a5251d41fd3c63e7892720e50d02aaef985fb9a.de-de.xlf", [Type]::Missing, [Type]::Missing, "8e2b4336-5735-4c5d-bffc-e21e7f60d885.2a5251d41fd3c63e7892720e50d02aaef985fb9a.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/d641bdf9e317d10099a875bca9a537f3b06a36e7/e2e/e0390748-aa6c-4338-96c4-c0d594e4314f.md", [Type]::Missing, [Type]::Missing, "e0390748-aa6c-4338-96c4-c0d594e4314f.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aa4e0fd75e94cdda9e4e2524ee477d04f655959d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e0390748-aa6c-4338-96c4-c0d594e4314f.6f38887ef8fb072b659eafcafeb8544e61d5db31.de-de.xlf", [Type]::Missing, [Type]::Missing, "e0390748-aa6c-4338-96c4-c0d594e4314f.6f38887ef8fb072b659eafcafeb8544e61d5db31.de-de.xlf") | Out-Null

Write-Output "Handback report generated"
